$d = $word.ActiveDocument

# 1. "What things you need to watch out when apply for bankruptcy?"
#    -> "What things you need consider seriously if you apply for bankruptcy?"
$d.Content.Find.Execute(
    "hings you need to watch out when apply",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "hings you need consider seriously if you apply",
    2)

# 2. "Your credit rating will be affected"
#    -> "Your credit rating will be affected."
$d.Content.Find.Execute(
    "Your credit rating will be affected",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Your credit rating will be affected.",
    2)

# 3. "The cost of bankruptcy is high (680 pounds)"
#    -> "The cost of bankruptcy is high, £ 680 pounds. "
$d.Content.Find.Execute(
    "The cost of bankruptcy is high (680 pounds)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The cost of bankruptcy is high, £ 680 pounds. ",
    2)

# 4. "Criminal fines, compensation orders and victim surcharges from a magistrates' court   or crown court"
#    -> "...magistrates' court or crown court" (collapse triple space to single space)
$d.Content.Find.Execute(
    "Criminal fines, compensation orders and victim surcharges from a magistrates’ court   or crown court",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Criminal fines, compensation orders and victim surcharges from a magistrates’ court or crown court",
    2)
